$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format while writing so that numeric-looking
# strings (e.g. "233.85") are not silently converted to floating point numbers.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range('D2').Value = '41.711.87'
$ws.Range('E2').Value = '  +5.69%  '
$ws.Range('D3').Value = '2.259.76'
$ws.Range('E3').Value = '  +4.42%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '233.85'
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('D6').Value = '0.645'
$ws.Range('E6').Value = '  +3.45%  '
$ws.Range('D7').Value = '64.35'
$ws.Range('E7').Value = '  +0.90%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +4.28%  '
$ws.Range('D10').Value = '60.20'
$ws.Range('E10').Value = '  +3.61%  '
$ws.Range('D11').Value = '0.0900'
$ws.Range('E11').Value = '  +5.36%  '
$ws.Range('E12').Value = '  +2.51%  '
$ws.Range('D13').Value = '2.596.77'
$ws.Range('E13').Value = '  +4.43%  '
$ws.Range('D14').Value = '16.13'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').Value = '22.77'
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('D16').Value = '0.828'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('E17').Value = '  +2.97%  '
$ws.Range('D18').Value = '2.264.15'
$ws.Range('E18').Value = '  +4.25%  '
$ws.Range('D19').Value = '41.609.27'
$ws.Range('E19').Value = '  +5.43%  '
$ws.Range('D20').Value = '0.0₃0939'
$ws.Range('E20').Value = '  +10.47%  '
$ws.Range('D21').Value = '74.74'
$ws.Range('E21').Value = '  +3.95%  '
$ws.Range('D22').Value = '6.18'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').Value = '252.48'
$ws.Range('E23').Value = '  +9.99%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.42'
$ws.Range('E25').Value = '  +3.52%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  +2.78%  '
$ws.Range('D27').Value = '9.86'
$ws.Range('E27').Value = '  +3.65%  '
$ws.Range('E28').Value = '  +4.93%  '
$ws.Range('D29').Value = '171.52'
$ws.Range('E30').Value = '  +3.40%  '
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('E32').Value = '  +6.55%  '
$ws.Range('E33').Value = '  +2.91%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.79'
$ws.Range('E34').Value = '  +3.50%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '5.06'
$ws.Range('E35').Value = '  +7.12%  '
$ws.Range('D36').Value = '0.0641'
$ws.Range('E36').Value = '  +3.31%  '
$ws.Range('D37').Value = '6.91'
$ws.Range('E37').Value = '  -2.51%  '
$ws.Range('D38').Value = '3.84'
$ws.Range('E38').Value = '  +7.09%  '
$ws.Range('D39').Value = '2.49'
$ws.Range('E39').Value = '  +1.71%  '
$ws.Range('D40').Value = '0.000260'
$ws.Range('E40').Value = '  +59.31%  '
$ws.Range('D41').Value = '5.10'
$ws.Range('E41').Value = '  +19.15%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').Value = '0.0242'
$ws.Range('E43').Value = '  +5.19%  '
$ws.Range('E44').Value = '  +13.17%  '
$ws.Range('D45').Value = '102.70'
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('E46').Value = '  +6.26%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '17.64'
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = '1.24'
$ws.Range('E48').Value = '  +2.32%  '
$ws.Range('D49').Value = '1.509.40'
$ws.Range('E49').Value = '  -0.99%  '
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('D51').Value = '2.80'
$ws.Range('E51').Value = '  -1.00%  '

# Restore the original (default) cell style now that the text values are set.
$priceVolumeRange.Style = "Normal"
